$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header "scraped_at" in column I, matching the style of the other
# header cells (bold, bordered, centered - same as H1).
$ws.Range("I1").Value = "scraped_at"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The data rows for "Medisch Pedicure Landgraaf" and "Hallux Medisch
# Pedicure" swap places: row 2 becomes the Hallux record, row 3 becomes
# the Landgraaf record. Column A ("organisatietype") is identical for
# both records, so it does not need to change.
$ws.Range("B2").Value = "Hallux Medisch Pedicure"
$ws.Range("C2").Value = "https://www.zorgkaartnederland.nl/zorginstelling/medisch-pedicurepraktijk-hallux-medisch-pedicure-maastricht-10048434"
$ws.Range("D2").Value = "Victor de Steursstraat"

# E2 ("huisnummer") looks numeric ("15"); force it to stay plain text
# without leaving any number-format styling behind.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "15"
$ws.Range("E2").Style = "Normal"

$ws.Range("G2").Value = "6217KP"
$ws.Range("H2").Value = "Maastricht"

$ws.Range("B3").Value = "Medisch Pedicure Landgraaf"
$ws.Range("C3").Value = "https://www.zorgkaartnederland.nl/zorginstelling/medisch-pedicurepraktijk-medisch-pedicure-landgraaf-landgraaf-10070714"
$ws.Range("D3").Value = "Kerkstraat"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "25"
$ws.Range("E3").Style = "Normal"

$ws.Range("G3").Value = "6374HH"
$ws.Range("H3").Value = "Landgraaf"

# Populate the new scraped_at column for every data row with the scrape
# date, keeping it as plain text (not auto-converted to an Excel date).
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2025-05-27"
$ws.Range("I2").Style = "Normal"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2025-05-27"
$ws.Range("I3").Style = "Normal"

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "2025-05-27"
$ws.Range("I4").Style = "Normal"
